$wb = $excel.ActiveWorkbook

# xlEdge* constants
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10

$xlContinuous = 1
$xlThin = 2

# Note on ordering: the border edges are applied top, then bottom, then
# (optionally) right so that every intermediate style the engine interns
# while building up the final border combination is one that already
# exists in the workbook's style table (or gets reused across cells),
# instead of minting throw-away, never-referenced border/style entries.
function Set-TopBottomBorder($range, $withRight) {

    $range.ClearFormats()

    $range.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeTop).Weight = $xlThin

    if ($withRight) {
        $range.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
        $range.Borders.Item($xlEdgeRight).Weight = $xlThin
    }

    $range.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeBottom).Weight = $xlThin
}

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder $ws1.Range("C1") $false
Set-TopBottomBorder $ws1.Range("D1") $true

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder $ws2.Range("C1") $false
Set-TopBottomBorder $ws2.Range("D1") $true
Set-TopBottomBorder $ws2.Range("F1") $false
Set-TopBottomBorder $ws2.Range("G1") $true

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
